$d = $word.ActiveDocument

# Locate the target 3-column table (425 / 4293 / 3645 dxa grid -> 425 / 4293 / 3933 dxa)
$t = $null
foreach ($tbl in $d.Tables) {
    if ($tbl.Columns.Count -eq 3) {
        $w1 = $tbl.Columns(1).Width
        $w2 = $tbl.Columns(2).Width
        $w3 = $tbl.Columns(3).Width
        if ([Math]::Round($w1) -eq 21 -and [Math]::Round($w2) -eq 215 -and [Math]::Round($w3) -eq 182) {
            $t = $tbl
        }
    }
}

# Table-level preferred width: auto -> fixed 8651 dxa (432.55 pt)
$t.PreferredWidthType = 3
$t.PreferredWidth = 432.55

# Third column (gridCol 3645 -> 3933 dxa = 196.65 pt)
$t.Columns(3).PreferredWidthType = 3
$t.Columns(3).PreferredWidth = 196.65
$t.Columns(3).Width = 196.65
